# convert the code to functions
# Update the Matches_A (column B) and Matches_B (column H) counts on the
# "Validation" sheet for the 4 teams listed in rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7
$ws.Range("H2").Value = 12

$ws.Range("B3").Value = 10
$ws.Range("H3").Value = 7

$ws.Range("B4").Value = 14
$ws.Range("H4").Value = 8

$ws.Range("B5").Value = 11
$ws.Range("H5").Value = 9
